# --- Corrected TOTAL row sums + new projected SMB / dH.dt / Snow precipitation columns (V-Y) ---
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename T/U headers and add new V/W/X/Y headers (row 1)
$ws.Range("T1").Value = "SMB.A1B.2001.2030.Gt.year"
$ws.Range("U1").Value = "SMB.A1B.2071.2100.Gt.year"
$ws.Range("V1").Value = "dH.dt.proj.2001.2030.Gt.year"
$ws.Range("W1").Value = "dH.dt.proj.2071.2100.Gt.year"
$ws.Range("X1").Value = "Snow.prec.2001.2030.mm.year"
$ws.Range("Y1").Value = "Snow.prec.2071.2100.mm.year"

# Row 2: recomputed SMB.A1B projections (T,U) + new dH.dt.proj (V,W) and Snow.prec (X,Y)
$ws.Range("T2").Value = 63.6212579037088
$ws.Range("U2").Value = 76.738035258348
$ws.Range("V2").Value = 184.521257903709
$ws.Range("W2").Value = 197.638035258348
$ws.Range("X2").Value = -543.831522495685
$ws.Range("Y2").Value = -582.490033065415

# Row 3: recomputed SMB.A1B projections (T,U) + new dH.dt.proj (V,W) and Snow.prec (X,Y)
$ws.Range("T3").Value = 13.9951919417901
$ws.Range("U3").Value = 16.0681775998268
$ws.Range("V3").Value = 70.7951919417901
$ws.Range("W3").Value = 72.8681775998268
$ws.Range("X3").Value = -679.069074473327
$ws.Range("Y3").Value = -698.953199561326

# Row 4: recomputed SMB.A1B projections (T,U) + new dH.dt.proj (V,W) and Snow.prec (X,Y)
$ws.Range("T4").Value = 24.0244524734754
$ws.Range("U4").Value = 26.5831133060038
$ws.Range("V4").Value = 37.1244524734754
$ws.Range("W4").Value = 39.6831133060038
$ws.Range("X4").Value = -449.995599656621
$ws.Range("Y4").Value = -481.009824485234

# Row 5: recomputed SMB.A1B projections (T,U) + new dH.dt.proj (V,W) and Snow.prec (X,Y)
$ws.Range("T5").Value = 8.37283269330014
$ws.Range("U5").Value = 8.76643606072401
$ws.Range("V5").Value = 13.8728326933001
$ws.Range("W5").Value = 14.266436060724
$ws.Range("X5").Value = -1082.72387151232
$ws.Range("Y5").Value = -1113.44317529396

# Row 6: recomputed SMB.A1B projections (T,U) + new dH.dt.proj (V,W) and Snow.prec (X,Y)
$ws.Range("T6").Value = 33.2852990878104
$ws.Range("U6").Value = 36.575581021163
$ws.Range("V6").Value = 119.98529908781
$ws.Range("W6").Value = 123.275581021163
$ws.Range("X6").Value = -1059.30413948028
$ws.Range("Y6").Value = -1088.35277542615

# Row 7: recomputed SMB.A1B projections (T,U) + new dH.dt.proj (V,W) and Snow.prec (X,Y)
$ws.Range("T7").Value = 4.99397571540916
$ws.Range("U7").Value = 5.27235417020563
$ws.Range("V7").Value = 25.1939757154092
$ws.Range("W7").Value = 25.4723541702056
$ws.Range("X7").Value = -3676.39247335608
$ws.Range("Y7").Value = -3717.01442471139

# Row 8: recomputed SMB.A1B projections (T,U) + new dH.dt.proj (V,W) and Snow.prec (X,Y)
$ws.Range("T8").Value = 8.4736795340646
$ws.Range("U8").Value = 8.29847288133104
$ws.Range("V8").Value = 49.5736795340646
$ws.Range("W8").Value = 49.398472881331
$ws.Range("X8").Value = -817.319212814729
$ws.Range("Y8").Value = -814.430587946896

# Row 9: recomputed SMB.A1B projections (T,U) + new dH.dt.proj (V,W) and Snow.prec (X,Y)
$ws.Range("T9").Value = 25.5136379510608
$ws.Range("U9").Value = 25.4622336872447
$ws.Range("V9").Value = 168.713637951061
$ws.Range("W9").Value = 168.662233687245
$ws.Range("X9").Value = -3537.84235296401
$ws.Range("Y9").Value = -3536.76443072927

# Row 10: recomputed SMB.A1B projections (T,U) + new dH.dt.proj (V,W) and Snow.prec (X,Y)
$ws.Range("T10").Value = 14.1167246460745
$ws.Range("U10").Value = 14.3732099231156
$ws.Range("V10").Value = 189.316724646074
$ws.Range("W10").Value = 189.573209923116
$ws.Range("X10").Value = -12890.8956403666
$ws.Range("Y10").Value = -12908.3601562238

# Row 11: recomputed SMB.A1B projections (T,U) + new dH.dt.proj (V,W) and Snow.prec (X,Y)
$ws.Range("T11").Value = 8.00980824200625
$ws.Range("U11").Value = 8.21281007405133
$ws.Range("V11").Value = 54.5098082420062
$ws.Range("W11").Value = 54.7128100740513
$ws.Range("X11").Value = -4408.4488078712
$ws.Range("Y11").Value = -4424.86646211247

# Row 12: recomputed SMB.A1B projections (T,U) + new dH.dt.proj (V,W) and Snow.prec (X,Y)
$ws.Range("T12").Value = 2.31109950611006
$ws.Range("U12").Value = 2.53312971432776
$ws.Range("V12").Value = 26.4110995061101
$ws.Range("W12").Value = 26.6331297143278
$ws.Range("X12").Value = -3505.95044352344
$ws.Range("Y12").Value = -3535.42391950638

# Row 13: recomputed SMB.A1B projections (T,U) + new dH.dt.proj (V,W) and Snow.prec (X,Y)
$ws.Range("T13").Value = 0.59913603139006
$ws.Range("U13").Value = 0.806042003053048
$ws.Range("V13").Value = 16.9991360313901
$ws.Range("W13").Value = 17.206042003053
$ws.Range("X13").Value = -3932.22171441322
$ws.Range("Y13").Value = -3980.08297942766

# Row 14: recomputed SMB.A1B projections (T,U) + new dH.dt.proj (V,W) and Snow.prec (X,Y)
$ws.Range("T14").Value = 28.8785200436443
$ws.Range("U14").Value = 33.0135691035203
$ws.Range("V14").Value = 88.2785200436443
$ws.Range("W14").Value = 92.4135691035203
$ws.Range("X14").Value = -453.398827879895
$ws.Range("Y14").Value = -474.636456195897

# Row 15: recomputed SMB.A1B projections (T,U) + new dH.dt.proj (V,W) and Snow.prec (X,Y)
$ws.Range("T15").Value = 32.9952297584975
$ws.Range("U15").Value = 34.7925819110992
$ws.Range("V15").Value = 4.19522975849748
$ws.Range("W15").Value = 5.99258191109923
$ws.Range("X15").Value = -13.7052155106575
$ws.Range("Y15").Value = -19.576907889378

# Row 16: recomputed SMB.A1B projections (T,U) + new dH.dt.proj (V,W) and Snow.prec (X,Y)
$ws.Range("T16").Value = 48.8478199407002
$ws.Range("U16").Value = 49.0638085134821
$ws.Range("V16").Value = 290.6478199407
$ws.Range("W16").Value = 290.863808513482
$ws.Range("X16").Value = -5282.75778136883
$ws.Range("Y16").Value = -5286.68354731396

# Row 17: recomputed SMB.A1B projections (T,U) + new dH.dt.proj (V,W) and Snow.prec (X,Y)
$ws.Range("T17").Value = 14.1744999406479
$ws.Range("U17").Value = 14.6414118713883
$ws.Range("V17").Value = 292.974499940648
$ws.Range("W17").Value = 293.441411871388
$ws.Range("X17").Value = -20162.4214038016
$ws.Range("Y17").Value = -20194.5541495114

# Row 18: recomputed SMB.A1B projections (T,U) + new dH.dt.proj (V,W) and Snow.prec (X,Y)
$ws.Range("T18").Value = 4.43619680759275
$ws.Range("U18").Value = 5.02303790497627
$ws.Range("V18").Value = 169.036196807593
$ws.Range("W18").Value = 169.623037904976
$ws.Range("X18").Value = -27048.3881344758
$ws.Range("Y18").Value = -27142.2917248018

# Row 19: recomputed SMB.A1B projections (T,U) + new dH.dt.proj (V,W) and Snow.prec (X,Y)
$ws.Range("T19").Value = 30.5026717992868
$ws.Range("U19").Value = 34.1534518765788
$ws.Range("V19").Value = 163.302671799287
$ws.Range("W19").Value = 166.953451876579
$ws.Range("X19").Value = -4532.08127637991
$ws.Range("Y19").Value = -4633.4000842731

# Row 20: recomputed SMB.A1B projections (T,U) + new dH.dt.proj (V,W) and Snow.prec (X,Y)
$ws.Range("T20").Value = 29.7959628954875
$ws.Range("U20").Value = 35.8009116933811
$ws.Range("V20").Value = 266.595962895487
$ws.Range("W20").Value = 272.600911693381
$ws.Range("X20").Value = -5451.24036807855
$ws.Range("Y20").Value = -5574.02699597717

# Row 21: recomputed SMB.A1B projections (T,U) + new dH.dt.proj (V,W) and Snow.prec (X,Y)
$ws.Range("T21").Value = 0.202688512456996
$ws.Range("U21").Value = 0.231335777956005
$ws.Range("V21").Value = 12.902688512457
$ws.Range("W21").Value = 12.931335777956
$ws.Range("X21").Value = -46605.3404820553
$ws.Range("Y21").Value = -46708.8162469063

# Row 22: recomputed SMB.A1B projections (T,U) + new dH.dt.proj (V,W) and Snow.prec (X,Y)
$ws.Range("T22").Value = 38.019217080966
$ws.Range("U22").Value = 27.9422037238043
$ws.Range("V22").Value = 93.819217080966
$ws.Range("W22").Value = 83.7422037238043
$ws.Range("X22").Value = -1206.20808366656
$ws.Range("Y22").Value = -1076.65067156265

# Row 23 (TOTAL): corrected column sums (B-N) + recomputed/new projection columns (T-Y)
$ws.Range("B23").Value = 1545836.915
$ws.Range("C23").Value = 1681.6
$ws.Range("D23").Value = 1088
$ws.Range("E23").Value = 429.1
$ws.Range("F23").Value = -287.4
$ws.Range("G23").Value = 1023.1
$ws.Range("H23").Value = 1310
$ws.Range("I23").Value = 135.9
$ws.Range("J23").Value = 128.1
$ws.Range("K23").Value = 79.4
$ws.Range("L23").Value = 333.9
$ws.Range("M23").Value = 203.1
$ws.Range("N23").Value = 423
$ws.Range("T23").Value = 435.16990250548
$ws.Range("U23").Value = 464.351908075581
$ws.Range("V23").Value = 2338.76990250548
$ws.Range("W23").Value = 2367.95190807558
$ws.Range("X23").Value = -1512.94737485647
$ws.Range("Y23").Value = -1531.82517838603
